$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4597200751304626
$ws.Range("B1").Value = 0.4164294004440308
$ws.Range("C1").Value = 3.45040225982666
$ws.Range("D1").Value = 1.626733660697937
$ws.Range("E1").Value = 1.145544528961182
